$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The account-statement table (rows 16-25 in the original) is being replaced:
#   - the "EDGARDO LUIS VERGARA CONTO" record moves from the bottom to the top
#   - a brand-new worker "EFRAIN ISAAC MARTINEZ MARTINEZ" is inserted right
#     after it
#   - the "WLADIMIR ALEXANDER GARCIA PERAZA" rows are kept but re-ordered
#     (periods now descend: 2309..2302)
#   - the "ALEXANDER MOISES ANILLO MONTES" record moves from the middle to
#     the new last row, keeping the special "closing" row formatting
#
# That grows the table from 10 data rows to 11, so we first insert one row
# right after the current last data row (25) to make room; this naturally
# pushes the two signature/footer rows (30, 31) down to (31, 32) and keeps
# the existing merged cells in sync.
# ---------------------------------------------------------------------------

$ws.Rows("26:26").Insert()

# Row 25 currently still carries the special "last row" border/format. Move
# that formatting down onto the freshly inserted row 26 (the new last row),
# then restore row 25 to the regular data-row formatting (copied from row 24).
$ws.Range("B25:J25").Copy()
$ws.Range("B26:J26").PasteSpecial(-4122)

$ws.Range("B24:J24").Copy()
$ws.Range("B25:J25").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Rewrite the data rows 16-26 with the new content.
# Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora, F=Valor Mora,
#          G=Salario Basico
# ---------------------------------------------------------------------------

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143375685"
$ws.Range("D16").Value = "EDGARDO LUIS VERGARA CONTO"
$ws.Range("E16").Value = "2311"
$ws.Range("F16").Value = 20107
$ws.Range("G16").Value = 1160000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1065569164"
$ws.Range("D17").Value = "EFRAIN ISAAC MARTINEZ MARTINEZ"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 1898
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "PPT"
$ws.Range("C18").Value = "1287309"
$ws.Range("D18").Value = "WLADIMIR ALEXANDER GARCIA PERAZA"
$ws.Range("E18").Value = "2309"
$ws.Range("F18").Value = 46400
$ws.Range("G18").Value = 1300000

$ws.Range("B19").Value = "PPT"
$ws.Range("C19").Value = "1287309"
$ws.Range("D19").Value = "WLADIMIR ALEXANDER GARCIA PERAZA"
$ws.Range("E19").Value = "2308"
$ws.Range("F19").Value = 46400
$ws.Range("G19").Value = 1300000

$ws.Range("B20").Value = "PPT"
$ws.Range("C20").Value = "1287309"
$ws.Range("D20").Value = "WLADIMIR ALEXANDER GARCIA PERAZA"
$ws.Range("E20").Value = "2307"
$ws.Range("F20").Value = 46400
$ws.Range("G20").Value = 1300000

$ws.Range("B21").Value = "PPT"
$ws.Range("C21").Value = "1287309"
$ws.Range("D21").Value = "WLADIMIR ALEXANDER GARCIA PERAZA"
$ws.Range("E21").Value = "2306"
$ws.Range("F21").Value = 46400
$ws.Range("G21").Value = 1300000

$ws.Range("B22").Value = "PPT"
$ws.Range("C22").Value = "1287309"
$ws.Range("D22").Value = "WLADIMIR ALEXANDER GARCIA PERAZA"
$ws.Range("E22").Value = "2305"
$ws.Range("F22").Value = 46400
$ws.Range("G22").Value = 1300000

$ws.Range("B23").Value = "PPT"
$ws.Range("C23").Value = "1287309"
$ws.Range("D23").Value = "WLADIMIR ALEXANDER GARCIA PERAZA"
$ws.Range("E23").Value = "2304"
$ws.Range("F23").Value = 46400
$ws.Range("G23").Value = 1300000

$ws.Range("B24").Value = "PPT"
$ws.Range("C24").Value = "1287309"
$ws.Range("D24").Value = "WLADIMIR ALEXANDER GARCIA PERAZA"
$ws.Range("E24").Value = "2303"
$ws.Range("F24").Value = 46400
$ws.Range("G24").Value = 1300000

$ws.Range("B25").Value = "PPT"
$ws.Range("C25").Value = "1287309"
$ws.Range("D25").Value = "WLADIMIR ALEXANDER GARCIA PERAZA"
$ws.Range("E25").Value = "2302"
$ws.Range("F25").Value = 46400
$ws.Range("G25").Value = 1300000

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "20429869"
$ws.Range("D26").Value = "ALEXANDER MOISES ANILLO MONTES"
$ws.Range("E26").Value = "2304"
$ws.Range("F26").Value = 18372
$ws.Range("G26").Value = 1531000

# ---------------------------------------------------------------------------
# Update the summary figures above the table:
#   - VALOR MORA total (E11): 409679 -> 411577
#   - Cant. Trabajadores (C13): 3 -> 4
#   - Cant. Periodos (F13): 9 -> 10
# ---------------------------------------------------------------------------

$ws.Range("E11").Value = 411577
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 10
